$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A (Stage) cells A3:A17 -> left-align; Column B (Surface area) B3:B17 -> scale by 1,000,000
for ($r = 3; $r -le 17; $r++) {
    $aCell = $ws.Cells.Item($r, 1)
    $aCell.HorizontalAlignment = -4131  # xlLeft

    $bCell = $ws.Cells.Item($r, 2)
    $bCell.Value2 = $bCell.Value2 * 1000000
}

# Update the active selection to G9
$ws.Range("G9").Select()
